$wb = $excel.ActiveWorkbook

# Style donors (kept on a sheet that is never touched) so the rebuilt
# "总计" sheet and the new "2022-Q1" sheet both end up with the same
# bold/centered/bordered header look ("s=2") used throughout the workbook.
$headerStyleSrc = $wb.Worksheets.Item(5).Range("B1:H1")
$indexStyleSrc5 = $wb.Worksheets.Item(5).Range("A2:A6")
$indexStyleSrc7 = $wb.Worksheets.Item(5).Range("A2:A7")

# ------------------------------------------------------------------
# The "总计" sheet needs a new row inserted above its existing data.
# Recreate it (delete + add) so the new "2022-Q1" sheet can slot in
# with sheetId 6 and "总计" naturally becomes sheetId 7, matching the
# order/numbering a user gets by inserting a sheet before it in Excel.
# ------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Delete()

# ------------------------------------------------------------------
# 1) Brand-new "2022-Q1" sheet, placed where "总计" used to be.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$q1 = $wb.Worksheets.Add($null, $lastSheet)
$q1.Name = "2022-Q1"

$headerStyleSrc.Copy()
$q1.Range("B1").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$indexStyleSrc5.Copy()
$q1.Range("A2:A6").PasteSpecial(-4122)

# Columns B-G hold text-formatted values (fund codes with leading
# zeros, decimal strings) in the source data, so force Text format
# before assigning to avoid Excel's automatic numeric coercion.
$q1.Range("B2:G6").NumberFormat = "@"

$fundRows = @(
    @(0, "377020", "上投摩根内需动力混合",         "21.62", "90.74", "3.68", "0.7956", 10),
    @(1, "008359", "华安医疗创新混合",             "6.73",  "92.10", "2.64", "0.1777", 9),
    @(2, "000073", "上投摩根成长动力混合",         "3.03",  "90.81", "4.14", "0.1254", 8),
    @(3, "005112", "银华中证全指医药卫生指数增强", "1.42",  "86.59", "7.71", "0.1095", 2),
    @(4, "000326", "南方中小盘成长股票",           "2.96",  "90.75", "1.90", "0.0562", 6)
)

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 2).Value = $row[1]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $q1.Cells.Item($r, 4).Value = $row[3]
    $q1.Cells.Item($r, 5).Value = $row[4]
    $q1.Cells.Item($r, 6).Value = $row[5]
    $q1.Cells.Item($r, 7).Value = $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Drop the temporary "@" number-format now that the text values are
# locked in, so the cells end up with plain (unstyled) formatting —
# matching the rest of the per-quarter sheets.
$q1.Range("B2:G6").ClearFormats()

# ------------------------------------------------------------------
# 2) Rebuild "总计" after "2022-Q1": original rows shifted down one,
#    with the new 2022-Q1 summary row on top.
# ------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$headerStyleSrc.Range("A1:C1").Copy()
$total.Range("B1").PasteSpecial(-4122)
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$indexStyleSrc7.Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q1", 5,  1.26),
    @(1, "2021-Q4", 16, 13.37),
    @(2, "2021-Q3", 17, 8.380000000000001),
    @(3, "2021-Q2", 15, 9.01),
    @(4, "2021-Q1", 19, 8.85),
    @(5, "2020-Q4", 30, 29.41)
)

$r = 2
foreach ($row in $summaryRows) {
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
